# Scheduled-runner update: refresh Market Board price/profit columns
# (currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ,
# LeveProfitNQ/HQ) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets.
# Mirrors a scheduled data-refresh commit; values below come from the
# latest market snapshot and replace the stale ones cell-by-cell.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (84 cells) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 663.6667
$ws.Range("I18").Value = 621.625
$ws.Range("K18").Value = 621.625
$ws.Range("M18").Value = -337.625
$ws.Range("H29").Value = 3483.8333
$ws.Range("I29").Value = 268.33334
$ws.Range("K29").Value = 805.0000200000001
$ws.Range("M29").Value = -524.0000200000001
$ws.Range("H31").Value = 508
$ws.Range("I31").Value = 508
$ws.Range("K31").Value = 1524
$ws.Range("M31").Value = -1294
$ws.Range("H32").Value = 1088.8572
$ws.Range("I32").Value = 1043.25
$ws.Range("J32").Value = 1149.6666
$ws.Range("K32").Value = 1043.25
$ws.Range("L32").Value = 1149.6666
$ws.Range("M32").Value = -717.25
$ws.Range("N32").Value = -1801.6666
$ws.Range("H34").Value = 7750
$ws.Range("J34").Value = 5000
$ws.Range("L34").Value = 5000
$ws.Range("N34").Value = -5406
$ws.Range("H36").Value = 7750
$ws.Range("J36").Value = 5000
$ws.Range("L36").Value = 5000
$ws.Range("N36").Value = -6430
$ws.Range("H38").Value = 10591.333
$ws.Range("J38").Value = 13349.333
$ws.Range("L38").Value = 40047.999
$ws.Range("N38").Value = -40791.999
$ws.Range("H40").Value = 1967
$ws.Range("J40").Value = 1967
$ws.Range("L40").Value = 1967
$ws.Range("N40").Value = -2317
$ws.Range("H41").Value = 2400.375
$ws.Range("I41").Value = 2933.2
$ws.Range("J41").Value = 1512.3334
$ws.Range("K41").Value = 2933.2
$ws.Range("L41").Value = 1512.3334
$ws.Range("M41").Value = -2493.2
$ws.Range("N41").Value = -2392.3334
$ws.Range("H42").Value = 665.6667
$ws.Range("I42").Value = 665.6667
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 1997.0001
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -1767.0001
$ws.Range("N42").ClearContents()
$ws.Range("H45").Value = 789
$ws.Range("I45").Value = 808.5
$ws.Range("K45").Value = 2425.5
$ws.Range("M45").Value = -2233.5
$ws.Range("H46").Value = 8993.5
$ws.Range("J46").Value = 5987
$ws.Range("L46").Value = 17961
$ws.Range("N46").Value = -18199
$ws.Range("H60").Value = 8993.5
$ws.Range("J60").Value = 5987
$ws.Range("L60").Value = 17961
$ws.Range("N60").Value = -18929
$ws.Range("H61").Value = 250
$ws.Range("I61").Value = 250
$ws.Range("K61").Value = 750
$ws.Range("M61").Value = -578
$ws.Range("H92").Value = 992.6667
$ws.Range("I92").Value = 989.5
$ws.Range("K92").Value = 989.5
$ws.Range("M92").Value = 258.5
$ws.Range("H100").Value = 1685
$ws.Range("I100").Value = 1890
$ws.Range("J100").Value = 660
$ws.Range("K100").Value = 1890
$ws.Range("L100").Value = 660
$ws.Range("M100").Value = -1349
$ws.Range("N100").Value = -1742
$ws.Range("H135").Value = 1159.6666
$ws.Range("I135").Value = 1128.8572
$ws.Range("K135").Value = 10159.7148
$ws.Range("M135").Value = -7624.7148
$ws.Range("H137").Value = 1487.375
$ws.Range("I137").Value = 1487.375
$ws.Range("K137").Value = 4462.125
$ws.Range("M137").Value = -1912.125

# ---- Sheet: ARM (19 cells) ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2431.75
$ws.Range("I45").Value = 2755.5
$ws.Range("J45").Value = 2323.8333
$ws.Range("K45").Value = 2755.5
$ws.Range("L45").Value = 2323.8333
$ws.Range("M45").Value = -2378.5
$ws.Range("N45").Value = -3077.8333
$ws.Range("H61").Value = 2999.3333
$ws.Range("J61").Value = 3000
$ws.Range("L61").Value = 3000
$ws.Range("N61").Value = -3424
$ws.Range("H94").Value = 30000
$ws.Range("J94").Value = 30000
$ws.Range("L94").Value = 30000
$ws.Range("N94").Value = -31802
$ws.Range("H136").Value = 2999.3333
$ws.Range("J136").Value = 3000
$ws.Range("L136").Value = 9000
$ws.Range("N136").Value = -14100

# ---- Sheet: BSM (16 cells) ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 685
$ws.Range("J64").Value = 595
$ws.Range("L64").Value = 595
$ws.Range("N64").Value = -1045
$ws.Range("H67").Value = 685
$ws.Range("J67").Value = 595
$ws.Range("L67").Value = 595
$ws.Range("N67").Value = -2155
$ws.Range("H107").Value = 1733.2858
$ws.Range("I107").Value = 1733.2858
$ws.Range("K107").Value = 1733.2858
$ws.Range("M107").Value = 186.7141999999999
$ws.Range("H134").Value = 1448.5
$ws.Range("I134").Value = 1448.5
$ws.Range("K134").Value = 4345.5
$ws.Range("M134").Value = -1810.5

# ---- Sheet: CRP (41 cells) ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2000
$ws.Range("I31").Value = 2000
$ws.Range("K31").Value = 2000
$ws.Range("M31").Value = -1705
$ws.Range("H34").Value = 2000
$ws.Range("I34").Value = 2000
$ws.Range("K34").Value = 2000
$ws.Range("M34").Value = -1798
$ws.Range("H41").Value = 8792.143
$ws.Range("J41").Value = 49999
$ws.Range("L41").Value = 49999
$ws.Range("N41").Value = -50855
$ws.Range("H58").Value = 4574.4165
$ws.Range("I58").Value = 2178.8
$ws.Range("J58").Value = 6285.5713
$ws.Range("K58").Value = 2178.8
$ws.Range("L58").Value = 6285.5713
$ws.Range("M58").Value = -1975.8
$ws.Range("N58").Value = -6691.5713
$ws.Range("H122").Value = 1946.375
$ws.Range("I122").Value = 1938.7142
$ws.Range("K122").Value = 5816.142599999999
$ws.Range("M122").Value = -3366.142599999999
$ws.Range("H132").Value = 5103.846
$ws.Range("J132").Value = 6285.7144
$ws.Range("L132").Value = 18857.1432
$ws.Range("N132").Value = -23917.1432
$ws.Range("H134").Value = 4669.077
$ws.Range("I134").Value = 1739.8
$ws.Range("J134").Value = 6499.875
$ws.Range("K134").Value = 5219.4
$ws.Range("L134").Value = 19499.625
$ws.Range("M134").Value = -2684.4
$ws.Range("N134").Value = -24569.625
$ws.Range("H136").Value = 4574.4165
$ws.Range("I136").Value = 2178.8
$ws.Range("J136").Value = 6285.5713
$ws.Range("K136").Value = 6536.400000000001
$ws.Range("L136").Value = 18856.7139
$ws.Range("M136").Value = -3986.400000000001
$ws.Range("N136").Value = -23956.7139

# ---- Sheet: CUL (27 cells) ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H14").Value = 498.5
$ws.Range("I14").Value = 498.5
$ws.Range("K14").Value = 1495.5
$ws.Range("M14").Value = -1322.5
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H129").Value = 5714.8335
$ws.Range("I129").Value = 1247
$ws.Range("J129").Value = 7948.75
$ws.Range("K129").Value = 3741
$ws.Range("L129").Value = 23846.25
$ws.Range("M129").Value = 1259
$ws.Range("N129").Value = -33846.25
$ws.Range("H136").Value = 2000
$ws.Range("J136").Value = 2000
$ws.Range("L136").Value = 6000
$ws.Range("N136").Value = -16200

# ---- Sheet: GSM (19 cells) ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("H29").Value = 425
$ws.Range("I29").Value = 425
$ws.Range("K29").Value = 425
$ws.Range("M29").Value = -135
$ws.Range("H102").Value = 1388.3572
$ws.Range("I102").Value = 1422.7273
$ws.Range("J102").Value = 1262.3334
$ws.Range("K102").Value = 1422.7273
$ws.Range("L102").Value = 1262.3334
$ws.Range("M102").Value = 199.2727
$ws.Range("N102").Value = -4506.3334
$ws.Range("H132").Value = 3500
$ws.Range("I132").Value = 3500
$ws.Range("K132").Value = 10500
$ws.Range("M132").Value = -7970

# ---- Sheet: LTW (18 cells) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 15248.75
$ws.Range("J31").Value = 25500
$ws.Range("L31").Value = 25500
$ws.Range("N31").Value = -25996
$ws.Range("H100").Value = 2495.5
$ws.Range("I100").Value = 2991
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 2991
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -2450
$ws.Range("N100").Value = -3082
$ws.Range("H136").Value = 31833.334
$ws.Range("I136").Value = 25250
$ws.Range("J136").Value = 45000
$ws.Range("K136").Value = 75750
$ws.Range("L136").Value = 135000
$ws.Range("M136").Value = -73200
$ws.Range("N136").Value = -140100
